$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" message on Hoja1!A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("1000 Bs = 3.34 = 13014.02 pesos", "1000 Bs = 3.34 = 13006.69 pesos")
$text = $text.Replace("13014.02 pesos = 3.33 = 970.23 Bs", "13006.69 pesos = 3.32 = 964.23 Bs")
$cell.Value = $text

# --- Update rate figures on tasas sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 299
$wsTasas.Range("O10").Value = 3889
$wsTasas.Range("N12").Value = 3911.88
$wsTasas.Range("O12").Value = 290
